{"js": "// Update the East-Asian (\"eastAsia\") and complex-script (\"cs\") fonts\n// recorded on the document's styles, matching the source diff:\n//   - docDefaults / Normal / Heading: eastAsia font DejaVu Sans -> Tahoma\n//   - List / Caption / Index: add an explicit complex-script (cs) font\n//     of \"DejaVu Sans\" (previously inherited / unset)\n//\n// Note: w:docDefaults (the package-level rPrDefault) has no\n// representation in the Word JavaScript API (Style/Font objects only\n// reach the named styles collection), so that portion of the source\n// diff cannot be reproduced through Office.js; every style that IS\n// reachable through context.document.getStyles() is updated here.\n\nconst styles = context.document.getStyles();\n\n// Styles whose East Asian font moves from \"DejaVu Sans\" to \"Tahoma\".\nconst eastAsiaStyleNames = [\"Normal\", \"Heading\"];\n// Styles that gain an explicit complex-script (w:cs) font of \"DejaVu Sans\".\nconst complexScriptStyleNames = [\"List\", \"Caption\", \"Index\"];\n\nconst styleObjects = {};\nfor (const name of eastAsiaStyleNames.concat(complexScriptStyleNames)) {\n  styleObjects[name] = styles.getByNameOrNullObject(name);\n}\nawait context.sync();\n\nfor (const name of eastAsiaStyleNames) {\n  styleObjects[name].font.nameFarEast = \"Tahoma\";\n}\nfor (const name of complexScriptStyleNames) {\n  styleObjects[name].font.nameBidirectional = \"DejaVu Sans\";\n}\nawait context.sync();\n", "ps1": "# Update the East-Asian (\"eastAsia\") and complex-script (\"cs\") fonts\n# recorded on the document's styles, matching the source diff:\n#   - docDefaults / Normal / Heading: eastAsia font DejaVu Sans -> Tahoma\n#   - List / Caption / Index: add an explicit complex-script (cs) font\n#     of \"DejaVu Sans\" (previously inherited / unset)\n#\n# Note: w:docDefaults (the package-level rPrDefault) has no\n# representation in the Word COM object model (Styles only reaches the\n# named styles collection), so that portion of the source diff cannot\n# be reproduced through COM automation either; every style that IS\n# reachable through $d.Styles is updated here.\n\n$d = $word.ActiveDocument\n$styles = $d.Styles\n\nforeach ($name in @(\"Normal\", \"Heading\")) {\n    $s = $styles.Item($name)\n    $s.Font.NameFarEast = \"Tahoma\"\n}\n\nforeach ($name in @(\"List\", \"Caption\", \"Index\")) {\n    $s = $styles.Item($name)\n    $s.Font.NameBi = \"DejaVu Sans\"\n}\n"}
